$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057161987598992
$ws.Range("D2").Value = 1.056392523886304
$ws.Range("E2").Value = 1.062507905838677
$ws.Range("F2").Value = 1.071210844270271
$ws.Range("I2").Value = 1.041693537971972
$ws.Range("J2").Value = 1.062159849014679
$ws.Range("K2").Value = 1.059130360032595
$ws.Range("L2").Value = 1.065229063986186
$ws.Range("M2").Value = 1.073908622668651
$ws.Range("N2").Value = 1.063668237831717

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058598364849935
$ws.Range("D3").Value = 1.057669166988266
$ws.Range("E3").Value = 1.06378723294479
$ws.Range("F3").Value = 1.072577604668824
$ws.Range("I3").Value = 1.041986010619345
$ws.Range("J3").Value = 1.063246244069036
$ws.Range("K3").Value = 1.060219454692441
$ws.Range("L3").Value = 1.066322053157196
$ws.Range("M3").Value = 1.075090533866597
$ws.Range("N3").Value = 1.064756175691663

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059526756863014
$ws.Range("D4").Value = 1.058494493375757
$ws.Range("E4").Value = 1.064614368254703
$ws.Range("F4").Value = 1.073461374708506
$ws.Range("I4").Value = 1.042173378908237
$ws.Range("J4").Value = 1.063947731201102
$ws.Range("K4").Value = 1.060922871893605
$ws.Range("L4").Value = 1.067028060664277
$ws.Range("M4").Value = 1.075854154256516
$ws.Range("N4").Value = 1.065458659015918

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059916809505704
$ws.Range("D5").Value = 1.058841286002138
$ws.Range("E5").Value = 1.064961938567528
$ws.Range("F5").Value = 1.073832768824097
$ws.Range("I5").Value = 1.042251699454659
$ws.Range("J5").Value = 1.064242285330271
$ws.Range("K5").Value = 1.061218281197161
$ws.Range("L5").Value = 1.067324575712942
$ws.Range("M5").Value = 1.076174907801562
$ws.Range("N5").Value = 1.065753631445736

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.059982286953653
$ws.Range("D6").Value = 1.058899503895692
$ws.Range("E6").Value = 1.065020288045626
$ws.Range("F6").Value = 1.073895119286382
$ws.Range("I6").Value = 1.042264823520805
$ws.Range("J6").Value = 1.064291721794112
$ws.Range("K6").Value = 1.061267863812665
$ws.Range("L6").Value = 1.067374344970523
$ws.Range("M6").Value = 1.07622874787398
$ws.Range("N6").Value = 1.065803138115027

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059531969714751
$ws.Range("D7").Value = 1.05849912791653
$ws.Range("E7").Value = 1.064619013119399
$ws.Range("F7").Value = 1.073466337850416
$ws.Range("I7").Value = 1.042174427194255
$ws.Range("J7").Value = 1.063951668422277
$ws.Range("K7").Value = 1.060926820369779
$ws.Range("L7").Value = 1.067032023849595
$ws.Range("M7").Value = 1.075858441246632
$ws.Range("N7").Value = 1.065462601828398

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.057647636784341
$ws.Range("D8").Value = 1.05682412865642
$ws.Range("E8").Value = 1.06294040269195
$ws.Range("F8").Value = 1.07167287747068
$ws.Range("I8").Value = 1.041792770715594
$ws.Range("J8").Value = 1.062527310811198
$ws.Range("K8").Value = 1.059498696173913
$ws.Range("L8").Value = 1.065598701710992
$ws.Range("M8").Value = 1.074308296746975
$ws.Range("N8").Value = 1.064036221466135

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054318981818133
$ws.Range("D9").Value = 1.053866652656911
$ws.Range("E9").Value = 1.059977119287495
$ws.Range("F9").Value = 1.068507656023554
$ws.Range("I9").Value = 1.041105780062409
$ws.Range("J9").Value = 1.060005871986236
$ws.Range("K9").Value = 1.056972027866646
$ws.Range("L9").Value = 1.063063414539128
$ws.Range("M9").Value = 1.071567712911009
$ws.Range("N9").Value = 1.061511201908829

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.052093995687927
$ws.Range("D10").Value = 1.051890755955904
$ws.Range("E10").Value = 1.057997728979847
$ws.Range("F10").Value = 1.066393918381038
$ws.Range("I10").Value = 1.040637981923954
$ws.Range("J10").Value = 1.058316901167422
$ws.Range("K10").Value = 1.055280523216043
$ws.Range("L10").Value = 1.061366530480681
$ws.Range("M10").Value = 1.06973432570529
$ws.Range("N10").Value = 1.059819832557725

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.051129076909486
$ws.Range("D11").Value = 1.051034102934246
$ws.Range("E11").Value = 1.057139653778298
$ws.Range("F11").Value = 1.065477730933281
$ws.Range("I11").Value = 1.040433076463136
$ws.Range("J11").Value = 1.057583600823438
$ws.Range("K11").Value = 1.054546352957779
$ws.Range("L11").Value = 1.060630118485963
$ws.Range("M11").Value = 1.068938888723061
$ws.Range("N11").Value = 1.059085490843142

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050770433026072
$ws.Range("D12").Value = 1.050715736766469
$ws.Range("E12").Value = 1.056820773038274
$ws.Range("F12").Value = 1.065137273514358
$ws.Range("I12").Value = 1.040356611389548
$ws.Range("J12").Value = 1.057310920221841
$ws.Range("K12").Value = 1.05427338349185
$ws.Range("L12").Value = 1.060356329556492
$ws.Range("M12").Value = 1.068643187281996
$ws.Range("N12").Value = 1.058812423003814

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050847373880705
$ws.Range("D13").Value = 1.050784035037446
$ws.Range("E13").Value = 1.056889181065999
$ws.Range("F13").Value = 1.06521030942089
$ws.Range("I13").Value = 1.040373029455839
$ws.Range("J13").Value = 1.057369424802161
$ws.Range("K13").Value = 1.054331948482277
$ws.Range("L13").Value = 1.060415069725575
$ws.Range("M13").Value = 1.068706627220978
$ws.Range("N13").Value = 1.058871010667348

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05109943602834
$ws.Range("D14").Value = 1.051007790127107
$ws.Range("E14").Value = 1.057113298156742
$ws.Range("F14").Value = 1.06544959159943
$ws.Range("I14").Value = 1.040426763068561
$ws.Range("J14").Value = 1.057561067113008
$ws.Range("K14").Value = 1.054523794670139
$ws.Range("L14").Value = 1.060607492191244
$ws.Range("M14").Value = 1.068914450890244
$ws.Range("N14").Value = 1.059062925132259

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.051254709138823
$ws.Range("D15").Value = 1.0511456307212
$ws.Range("E15").Value = 1.057251363607961
$ws.Range("F15").Value = 1.065597001916215
$ws.Range("I15").Value = 1.040459823155011
$ws.Range("J15").Value = 1.057679104335381
$ws.Range("K15").Value = 1.054641962051287
$ws.Range("L15").Value = 1.060726016389341
$ws.Range("M15").Value = 1.069042465859203
$ws.Range("N15").Value = 1.059181129981027

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.052158002286356
$ws.Range("D16").Value = 1.0519475860565
$ws.Range("E16").Value = 1.058054655381035
$ws.Range("F16").Value = 1.066454702754769
$ws.Range("I16").Value = 1.040651531262381
$ws.Range("J16").Value = 1.058365526072572
$ws.Range("K16").Value = 1.055329210652575
$ws.Range("L16").Value = 1.061415368544162
$ws.Range("M16").Value = 1.069787082761092
$ws.Range("N16").Value = 1.059868526515818

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052724211704222
$ws.Range("D17").Value = 1.052450339433227
$ws.Range("E17").Value = 1.058558271376235
$ws.Range("F17").Value = 1.066992464197499
$ws.Range("I17").Value = 1.040771155510771
$ws.Range("J17").Value = 1.058795570724989
$ws.Range("K17").Value = 1.055759834900216
$ws.Range("L17").Value = 1.061847335785163
$ws.Range("M17").Value = 1.070253737808083
$ws.Range("N17").Value = 1.060299181880969

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053054329260945
$ws.Range("D18").Value = 1.052743483447638
$ws.Range("E18").Value = 1.058851927275089
$ws.Range("F18").Value = 1.067306042472662
$ws.Range("I18").Value = 1.040840704088715
$ws.Range("J18").Value = 1.059046219257867
$ws.Range("K18").Value = 1.056010843419622
$ws.Range("L18").Value = 1.062099136096338
$ws.Range("M18").Value = 1.070525778839497
$ws.Range("N18").Value = 1.060550186363511

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.053166866802268
$ws.Range("D19").Value = 1.052843420581833
$ws.Range("E19").Value = 1.058952040427594
$ws.Range("F19").Value = 1.06741294965509
$ws.Range("I19").Value = 1.040864380034397
$ws.Range("J19").Value = 1.05913165198306
$ws.Range("K19").Value = 1.056096402605091
$ws.Range("L19").Value = 1.062184966693228
$ws.Range("M19").Value = 1.070618512349268
$ws.Range("N19").Value = 1.060635740412975

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052663477615082
$ws.Range("D20").Value = 1.052396409500658
$ws.Range("E20").Value = 1.058504247988918
$ws.Range("F20").Value = 1.066934776717762
$ws.Range("I20").Value = 1.040758344372157
$ws.Range("J20").Value = 1.058749450613771
$ws.Range("K20").Value = 1.055713650335929
$ws.Range("L20").Value = 1.061801006256316
$ws.Range("M20").Value = 1.070203685789884
$ws.Range("N20").Value = 1.060252996273902

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.051025216416009
$ws.Range("D21").Value = 1.050941904485716
$ws.Range("E21").Value = 1.057047305536529
$ws.Range("F21").Value = 1.065379132999012
$ws.Range("I21").Value = 1.040410949645928
$ws.Range("J21").Value = 1.057504641556511
$ws.Range("K21").Value = 1.054467308130228
$ws.Range("L21").Value = 1.060550835586548
$ws.Range("M21").Value = 1.068853258686112
$ws.Range("N21").Value = 1.059006419444999

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049993841888705
$ws.Range("D22").Value = 1.050026430038811
$ws.Range("E22").Value = 1.056130377814876
$ws.Range("F22").Value = 1.064400198583066
$ws.Range("I22").Value = 1.040190479480823
$ws.Range("J22").Value = 1.056720241304336
$ws.Range("K22").Value = 1.053682142207046
$ws.Range("L22").Value = 1.059763339056088
$ws.Range("M22").Value = 1.068002796222028
$ws.Range("N22").Value = 1.058220905254494

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050540721889181
$ws.Range("D23").Value = 1.050511834041703
$ws.Range("E23").Value = 1.056616544688155
$ws.Range("F23").Value = 1.064919231696583
$ws.Range("I23").Value = 1.040307549686804
$ws.Range("J23").Value = 1.057136233274169
$ws.Range("K23").Value = 1.054098521239137
$ws.Range("L23").Value = 1.060180946316111
$ws.Range("M23").Value = 1.068453776446177
$ws.Range("N23").Value = 1.058637487980639

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052690921177331
$ws.Range("D24").Value = 1.052420778436041
$ws.Range("E24").Value = 1.058528659125499
$ws.Range("F24").Value = 1.066960843481731
$ws.Range("I24").Value = 1.040764133873386
$ws.Range("J24").Value = 1.058770290891851
$ws.Range("K24").Value = 1.055734519669852
$ws.Range("L24").Value = 1.061821941067389
$ws.Range("M24").Value = 1.070226302609247
$ws.Range("N24").Value = 1.060273866147567

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055180529548909
$ws.Range("D25").Value = 1.05463195956675
$ws.Range("E25").Value = 1.060743861236948
$ws.Range("F25").Value = 1.069326553372685
$ws.Range("I25").Value = 1.041285105785789
$ws.Range("J25").Value = 1.060659117768729
$ws.Range("K25").Value = 1.057626458155933
$ws.Range("L25").Value = 1.063720008637847
$ws.Range("M25").Value = 1.072277316353351
$ws.Range("N25").Value = 1.062165375375262
